$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "70.288.12"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.603.47"
$ws.Range("E3").Value = "  +1.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "604.60"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6
Set-TextValue $ws.Range("D6") "195.63"
$ws.Range("E6").Value = "  -0.86%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
Set-TextValue $ws.Range("D8") "1.00"

# Row 9
Set-TextValue $ws.Range("D9") "0.207"
$ws.Range("E9").Value = "  -1.91%  "

# Row 10
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
Set-TextValue $ws.Range("D11") "53.76"
$ws.Range("E11").Value = "  -0.85%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000304"
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
Set-TextValue $ws.Range("D13") "9.55"
$ws.Range("E13").Value = "  -0.36%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.176.08"
$ws.Range("E14").Value = "  +1.65%  "

# Row 15
Set-TextValue $ws.Range("D15") "13.06"
$ws.Range("E15").Value = "  +2.63%  "

# Row 16
Set-TextValue $ws.Range("D16") "595.92"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17
Set-TextValue $ws.Range("D17") "70.383.75"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D18") "19.04"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D19") "3.600.39"
$ws.Range("E19").Value = "  +1.30%  "

# Row 20
$ws.Range("E20").Value = "  +1.32%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.996"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
Set-TextValue $ws.Range("D22") "17.80"
$ws.Range("E22").Value = "  -2.26%  "

# Row 23
$ws.Range("E23").Value = "  -2.70%  "

# Row 24
Set-TextValue $ws.Range("D24") "102.11"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25
Set-TextValue $ws.Range("D25") "4.62"
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.03"
$ws.Range("E26").Value = "  -2.99%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.77"
$ws.Range("E27").Value = "  -1.89%  "

# Row 28
$ws.Range("E28").Value = "  -0.80%  "

# Row 29
Set-TextValue $ws.Range("D29") "33.82"
$ws.Range("E29").Value = "  +0.37%  "

# Row 30
Set-TextValue $ws.Range("D30") "4.80"
$ws.Range("E30").Value = "  +6.40%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.20"
$ws.Range("E31").Value = "  +0.81%  "

# Row 32
Set-TextValue $ws.Range("D32") "12.29"
$ws.Range("E32").Value = "  -3.72%  "

# Row 33
$ws.Range("E33").Value = "  +0.95%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "63.16"
$ws.Range("E34").Value = "  -0.56%  "

# Row 35
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D35") "0.0₃0900"
$ws.Range("E35").Value = "  +7.98%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.890.36"
$ws.Range("E36").Value = "  +3.96%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D37") "529.33"
$ws.Range("E37").Value = "  +5.88%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D38") "3.10"
$ws.Range("E38").Value = "  -0.44%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.999"
$ws.Range("E39").Value = "  -0.14%  "

# Row 40
Set-TextValue $ws.Range("D40") "36.93"
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.391"
$ws.Range("E41").Value = "  -1.31%  "

# Row 42
Set-TextValue $ws.Range("D42") "3.53"
$ws.Range("E42").Value = "  -3.08%  "

# Row 43
$ws.Range("E43").Value = "  -1.97%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0454"
$ws.Range("E44").Value = "  -0.98%  "

# Row 45
Set-TextValue $ws.Range("D45") "3.40"
$ws.Range("E45").Value = "  +2.27%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.86"
$ws.Range("E46").Value = "  +0.65%  "

# Row 47
$ws.Range("E47").Value = "  +0.04%  "

# Row 48
Set-TextValue $ws.Range("D48") "8.63"
$ws.Range("E48").Value = "  -0.75%  "

# Row 49
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.000251"
$ws.Range("E50").Value = "  +1.85%  "

# Row 51
$ws.Range("E51").Value = "  +0.47%  "
